$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.738.48"
$ws.Range("E2").Value = "  +2.07%  "
$ws.Range("D3").Value = "3.081.86"
$ws.Range("E3").Value = "  +4.37%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.85"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.85"
$ws.Range("E6").Value = "  +5.18%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.079.21"
$ws.Range("E8").Value = "  +4.51%  "
$ws.Range("E9").Value = "  +0.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.59"
$ws.Range("E11").Value = "  +1.46%  "
$ws.Range("E12").Value = "  +4.87%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000249"
$ws.Range("E13").Value = "  +1.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.44"
$ws.Range("E14").Value = "  +6.26%  "
$ws.Range("E15").Value = "  -0.61%  "
$ws.Range("D16").Value = "3.592.97"
$ws.Range("E16").Value = "  +4.36%  "
$ws.Range("D17").Value = "66.671.36"
$ws.Range("E17").Value = "  +2.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.20"
$ws.Range("E18").Value = "  +3.20%  "
$ws.Range("D19").Value = "3.083.30"
$ws.Range("E19").Value = "  +4.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.20"
$ws.Range("E20").Value = "  +10.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "464.81"
$ws.Range("E21").Value = "  +4.18%  "
$ws.Range("E22").Value = "  +3.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.54"
$ws.Range("E23").Value = "  +3.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.18"
$ws.Range("E24").Value = "  +1.07%  "
$ws.Range("E25").Value = "  +5.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.89"
$ws.Range("E26").Value = "  +6.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.15"
$ws.Range("E27").Value = "  +1.19%  "
$ws.Range("E29").Value = "  -1.10%  "
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.66"
$ws.Range("E31").Value = "  +3.04%  "
$ws.Range("E32").Value = "  +1.09%  "
$ws.Range("E33").Value = "  +3.60%  "
$ws.Range("E34").Value = "  +3.13%  "
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("E36").Value = "  +2.81%  "
$ws.Range("E37").Value = "  +2.59%  "
$ws.Range("E38").Value = "  +7.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "46.94"
$ws.Range("E39").Value = "  +6.21%  "
$ws.Range("E40").Value = "  +6.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.14"
$ws.Range("E41").Value = "  +2.26%  "
$ws.Range("E42").Value = "  +1.60%  "
$ws.Range("E43").Value = "  +2.29%  "
$ws.Range("E44").Value = "  -1.37%  "
$ws.Range("E45").Value = "  +2.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "383.09"
$ws.Range("E46").Value = "  -0.59%  "
$ws.Range("D47").Value = "2.760.78"
$ws.Range("E47").Value = "  +1.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.62"
$ws.Range("E48").Value = "  +1.08%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("E50").Value = "  +5.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.24"
$ws.Range("E51").Value = "  +2.69%  "
